$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number (e.g. "571.69")
# need to be forced to Text format first, otherwise Excel auto-converts
# the assigned string into a numeric value - the source data keeps these
# as text cells (t="inlineStr") even when they look numeric.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range('D2').Value = '59.855.31'
$ws.Range('E2').Value = '  -5.20%  '
$ws.Range('D3').Value = '2.963.62'
$ws.Range('E3').Value = '  -6.89%  '
$ws.Range('E4').Value = '  -0.03%  '
Set-TextValue $ws.Range('D5') '571.69'
$ws.Range('E5').Value = '  -3.43%  '
Set-TextValue $ws.Range('D6') '124.73'
$ws.Range('E6').Value = '  -7.83%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '2.960.34'
$ws.Range('E8').Value = '  -6.90%  '
Set-TextValue $ws.Range('D9') '0.502'
$ws.Range('E9').Value = '  -2.61%  '
Set-TextValue $ws.Range('D10') '0.132'
$ws.Range('E10').Value = '  -6.51%  '
$ws.Range('E11').Value = '  -3.44%  '
Set-TextValue $ws.Range('D12') '0.438'
$ws.Range('E12').Value = '  -3.51%  '
$ws.Range('E13').Value = '  -6.90%  '
Set-TextValue $ws.Range('D14') '32.41'
$ws.Range('E14').Value = '  -6.52%  '
Set-TextValue $ws.Range('D15') '0.119'
$ws.Range('E15').Value = '  -0.66%  '
$ws.Range('D16').Value = '3.451.60'
$ws.Range('E16').Value = '  -6.87%  '
$ws.Range('D17').Value = '2.964.50'
$ws.Range('E17').Value = '  -6.93%  '
$ws.Range('D18').Value = '59.803.86'
$ws.Range('E18').Value = '  -5.18%  '
Set-TextValue $ws.Range('D19') '6.17'
$ws.Range('E19').Value = '  -6.10%  '
Set-TextValue $ws.Range('D20') '433.07'
$ws.Range('E20').Value = '  -6.25%  '
Set-TextValue $ws.Range('D21') '13.06'
$ws.Range('E21').Value = '  -7.03%  '
Set-TextValue $ws.Range('D22') '0.659'
$ws.Range('E22').Value = '  -5.66%  '
Set-TextValue $ws.Range('D23') '6.97'
$ws.Range('E23').Value = '  -8.84%  '
Set-TextValue $ws.Range('D24') '12.70'
$ws.Range('E24').Value = '  -5.14%  '
Set-TextValue $ws.Range('D25') '79.06'
$ws.Range('E25').Value = '  -4.33%  '
$ws.Range('E26').Value = '  +0.03%  '
Set-TextValue $ws.Range('D27') '0.998'
$ws.Range('E27').Value = '  -0.25%  '
Set-TextValue $ws.Range('D28') '2.51'
$ws.Range('E28').Value = '  -6.01%  '
Set-TextValue $ws.Range('D29') '7.21'
$ws.Range('E29').Value = '  -6.10%  '
$ws.Range('E30').Value = '  -7.82%  '
$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D31') '6.13'
$ws.Range('E31').Value = '  -9.36%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D32') '25.27'
$ws.Range('E32').Value = '  -9.27%  '
$ws.Range('E33').Value = '  -9.82%  '
Set-TextValue $ws.Range('D34') '2.17'
$ws.Range('E34').Value = '  -8.67%  '
Set-TextValue $ws.Range('D35') '0.943'
$ws.Range('E35').Value = '  -8.94%  '
Set-TextValue $ws.Range('D36') '5.57'
$ws.Range('E36').Value = '  -4.32%  '
Set-TextValue $ws.Range('D37') '49.56'
$ws.Range('E37').Value = '  -3.57%  '
$ws.Range('D38').Value = '0.0₃0650'
$ws.Range('E38').Value = '  -8.26%  '
Set-TextValue $ws.Range('D39') '7.93'
$ws.Range('E39').Value = '  -2.03%  '
$ws.Range('E40').Value = '  -8.52%  '
Set-TextValue $ws.Range('D41') '0.109'
$ws.Range('E41').Value = '  -3.18%  '
Set-TextValue $ws.Range('D42') '379.66'
$ws.Range('E42').Value = '  -6.61%  '
Set-TextValue $ws.Range('D43') '2.45'
$ws.Range('E43').Value = '  -7.51%  '
$ws.Range('D44').Value = '2.620.77'
$ws.Range('E44').Value = '  -6.85%  '
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('E46').Value = '  -6.68%  '
Set-TextValue $ws.Range('D47') '118.54'
$ws.Range('E47').Value = '  -4.62%  '
$ws.Range('E48').Value = '  -6.93%  '
$ws.Range('E49').Value = '  -4.44%  '
Set-TextValue $ws.Range('D50') '23.31'
$ws.Range('E50').Value = '  -7.67%  '
Set-TextValue $ws.Range('D51') '31.21'
$ws.Range('E51').Value = '  -11.25%  '
